$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed
$ws.Range("A2").Value = 4000
$ws.Range("B2").Value = 4000

$ws.Range("A4").Value = 4000
$ws.Range("B4").Value = 4000

$ws.Range("A8").Value = 6000
$ws.Range("B8").Value = 6000

# Add new rows 9-14
$ws.Range("A9").Value = 6000
$ws.Range("B9").Value = 6000

$ws.Range("A10").Value = 4000
$ws.Range("B10").Value = 4000

$ws.Range("A11").Value = 4000
$ws.Range("B11").Value = 4000

$ws.Range("A12").Value = 4000
$ws.Range("B12").Value = 4000

$ws.Range("A13").Value = 4000
$ws.Range("B13").Value = 4000

$ws.Range("A14").Value = 4000
$ws.Range("B14").Value = 4000
